# Weekly update: insert a new daily price record as row 114, pushing the
# existing records (old rows 114-148) down by one row to rows 115-149.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 114 - this shifts the
# existing rows 114:148 down to 115:149 and extends the used range to
# A1:R149, inheriting the row-above formatting (e.g. the date style on D).
$ws.Rows("114:114").Insert()

# Populate the new row 114 with the new weekly record.
$ws.Cells.Item(114, 1).Value = 4
$ws.Cells.Item(114, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(114, 3).Value = "Los Lagos"
$ws.Cells.Item(114, 4).Value = 44463
$ws.Cells.Item(114, 5).Value = 10
$ws.Cells.Item(114, 6).Value = 100112003
$ws.Cells.Item(114, 7).Value = "Ajo"
$ws.Cells.Item(114, 8).Value = "Chino"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 240
$ws.Cells.Item(114, 11).Value = 17000
$ws.Cells.Item(114, 12).Value = 17000
$ws.Cells.Item(114, 13).Value = 17000
$ws.Cells.Item(114, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(114, 15).Value = "China"
$ws.Cells.Item(114, 16).Value = 1700
$ws.Cells.Item(114, 17).Value = 10
$ws.Cells.Item(114, 18).Value = "Hortaliza"
